$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "Sheet1"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet1"

# Populate column A: header "WATER_USES" followed by the de-duplicated,
# alphabetically sorted list of WATER_USES letter-combination codes
# (built from the raw data's I/S/D/M/X/P/O beneficial-use letters).
$values = @("WATER_USES","D","DI","DIM","DIMO","DIMOP","DIMOS","DIMOSP","DIMOSPX","DIMOSX","DIMOX","DIMS","DIMSP","DIMX","DIO","DIOP","DIOS","DIOSP","DIOSPX","DIOSX","DIOX","DIP","DIS","DISP","DISPX","DISX","DIX","DM","DMO","DMOP","DMOS","DMOSX","DMS","DMX","DO","DOP","DOPX","DOS","DOSX","DOX","DP","DPX","DS","DSP","DSX","DX","I","IM","IMO","IMOP","IMOS","IMOSP","IMOSX","IMP","IMS","IMSP","IMX","IO","IOP","IOS","IOSP","IOSPX","IOSX","IOX","IP","IPX","IS","ISP","ISX","IX","M","MO","MOP","MOS","MOX","MP","MS","MSP","MSX","MX","O","OP","OPX","OS","OSP","OSX","OX","P","PX","S","SP","SX","X")

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Re-apply the A-Z sort on the data range (mirrors the author's manual
# Data > Sort step) so the sheet records the same sort state.
$sortRange = $ws.Range("A2:A317487")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortRange)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Make this the active/selected sheet with the same selection the author left.
$ws.Activate()
$ws.Range("A2:A317488").Select()
